$d = $word.ActiveDocument

$pairs = @(
    @("157×6=942", "391×5=1955"),
    @("632×2=1264", "207×9=1863"),
    @("924×7=6468", "690×3=2070"),
    @("660×4=2640", "737×4=2948"),
    @("406×2=812", "596×8=4768"),
    @("646×2=1292", "104×4=416"),
    @("615×5=3075", "622×7=4354"),
    @("453×3=1359", "640×9=5760"),
    @("293×6=1758", "174×5=870"),
    @("694×8=5552", "705×3=2115"),
    @("347×2=694", "483×6=2898"),
    @("360×7=2520", "196×7=1372"),
    @("823×3=2469", "952×6=5712"),
    @("788×6=4728", "213×5=1065"),
    @("878×3=2634", "580×8=4640"),
    @("886×3=2658", "914×3=2742"),
    @("730×2=1460", "510×7=3570"),
    @("907×4=3628", "101×8=808"),
    @("225×3=675", "424×7=2968"),
    @("309×9=2781", "796×2=1592"),
    @("900×6=5400", "245×9=2205"),
    @("841×2=1682", "934×9=8406"),
    @("881×3=2643", "697×7=4879"),
    @("415×4=1660", "668×6=4008"),
    @("842×8=6736", "782×8=6256")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done applying $($pairs.Count) replacements"
